# Update "想去人数" (attendance/wish-count) figures across the workbook.
# Sheet names: 展览 (Exhibition), 演出 (Show), 本地生活 (Local Life), 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param(
        [string]$SheetName,
        [string]$CellRef,
        [double]$NewValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $NewValue
}

# --- Sheet: 展览 ---
$sheet = "展览"
Set-CellValue $sheet "F3" 3824
Set-CellValue $sheet "F5" 1368
Set-CellValue $sheet "F10" 8575
Set-CellValue $sheet "F14" 110
Set-CellValue $sheet "F15" 289
Set-CellValue $sheet "F16" 330
Set-CellValue $sheet "F17" 85
Set-CellValue $sheet "F18" 354
Set-CellValue $sheet "F19" 10833
Set-CellValue $sheet "F23" 185
Set-CellValue $sheet "F28" 2673
Set-CellValue $sheet "F35" 2463
Set-CellValue $sheet "F36" 277
Set-CellValue $sheet "F37" 2579
Set-CellValue $sheet "F38" 3025
Set-CellValue $sheet "F39" 1242
Set-CellValue $sheet "F42" 340
Set-CellValue $sheet "F43" 314
Set-CellValue $sheet "F47" 88
Set-CellValue $sheet "F48" 96
Set-CellValue $sheet "F49" 81

# --- Sheet: 演出 ---
$sheet = "演出"
Set-CellValue $sheet "F9" 7
Set-CellValue $sheet "F11" 32
Set-CellValue $sheet "F16" 11

# --- Sheet: 本地生活 ---
$sheet = "本地生活"
Set-CellValue $sheet "F3" 30

# --- Sheet: 全部类型 ---
$sheet = "全部类型"
Set-CellValue $sheet "F3" 3824
Set-CellValue $sheet "F6" 1368
Set-CellValue $sheet "F11" 8575
Set-CellValue $sheet "F13" 110
Set-CellValue $sheet "F14" 289
Set-CellValue $sheet "F15" 330
Set-CellValue $sheet "F16" 85
Set-CellValue $sheet "F17" 354
Set-CellValue $sheet "F18" 10833
Set-CellValue $sheet "F20" 30
Set-CellValue $sheet "F23" 185
Set-CellValue $sheet "F29" 2673
Set-CellValue $sheet "F35" 2464
Set-CellValue $sheet "F36" 277
Set-CellValue $sheet "F37" 2579
Set-CellValue $sheet "F38" 3025
Set-CellValue $sheet "F40" 1242
Set-CellValue $sheet "F43" 340
Set-CellValue $sheet "F44" 314
Set-CellValue $sheet "F47" 88
Set-CellValue $sheet "F48" 96
Set-CellValue $sheet "F49" 81
